# Update the "想去人数" (want-to-go count) values for the first two events
# on the "展览" and "全部类型" sheets, as published by the gh-pages data
# generator (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 498
    $ws.Range("F3").Value = 3381
}
